$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fix grammatical errors / capitalization in the flow-of-events texts
$ws.Range("B6").Value = "The taxi driver accepts the request"
$ws.Range("B5").Value = "The system register the request and 10 minutes before the scheduled deparure send the request's and the user's basic information to the first taxi driver in the zone"
$ws.Range("B7").Value = "The system removes the taxi driver from the queue and send to the user the reminder of the incoming taxi and the expected waiting time"

# Update the saved selection/active cell
$ws.Range("B17").Select()
